# "Generate Report for Handoff"
#
# Updates the "Latest Handoff Datetime" (column D) for the row that tracks the
# "9015650a-e053-4466-aeaa-9d9227d4ceb1" file (row 4) on both the "zh-cn" and
# "de-de" localization-status sheets, reflecting the timestamp of the newly
# generated handoff report.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D4").Value = "2016-03-09 07:47:11"
$wsDe.Range("D4").Value = "2016-03-09 07:47:22"
